$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.431.79'
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("D3").Value = '1.862.32'
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.32'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4779'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3812'
$ws.Range("E8").Value = '  +3.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07329'
$ws.Range("E9").Value = '  +1.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9344'
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("E11").Value = '  +5.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07808'
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").Value = '1.903.84'
$ws.Range("E13").Value = '  +2.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.436'
$ws.Range("E14").Value = '  +1.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.553'
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.30'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008821'
$ws.Range("E18").Value = '  +1.84%  '
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = '27.515.68'
$ws.Range("E20").Value = '  +1.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.66'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.71'
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.943'
$ws.Range("E24").Value = '  +0.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.43'
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.47'
$ws.Range("E26").Value = '  +1.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.015'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.37'
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.943'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08888'
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.321'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.215'
$ws.Range("E32").Value = '  +3.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7548'
$ws.Range("E33").Value = '  +2.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.598'
$ws.Range("E34").Value = '  +2.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.713'
$ws.Range("E35").Value = '  -1.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02051'
$ws.Range("E36").Value = '  +4.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.125'
$ws.Range("E37").Value = '  +1.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5551'
$ws.Range("E38").Value = '  +5.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05275'
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.995'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.065'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.640'
$ws.Range("E42").Value = '  +4.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1524'
$ws.Range("E43").Value = '  +0.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4898'
$ws.Range("E44").Value = '  +3.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.72'
$ws.Range("E45").Value = '  +0.63%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.664'
$ws.Range("E47").Value = '  +3.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.01'
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.51'
$ws.Range("E49").Value = '  +2.98%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9152'
$ws.Range("E51").Value = '  +2.95%  '
